$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D (this pushes old D->E, E->F, F->G and copies
#    formatting from column C into the new column D, matching the diff).
$ws.Columns("D:D").Insert()

# 2. New column D should have the same width as column C.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# 3. Set the new header text in D20 ("Tỉ lệ điểm (%)").
$ws.Range("D20").Value = "Tỉ lệ điểm (%)"

# 4. Add the new "Điểm chữ:" row (row 26) below the "Tổng cộng:" row.
$ws.Range("B26:E26").Merge()
$ws.Range("B26").Value = "Điểm chữ:"
$ws.Range("B26:E26").HorizontalAlignment = -4152
$ws.Range("B26:E26").Font.Bold = $true
$ws.Range("B26:E26").Borders.LineStyle = 1
$ws.Range("F26").Borders.LineStyle = 1
$ws.Rows("26:26").RowHeight = $ws.Rows("25:25").RowHeight

# 5. Update the selection / view.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("F26").Select()

# 6. Update the print scale.
$ws.PageSetup.Zoom = 70
